$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that originally sits after
#    "en la estrategia A." (it moves, in the new document, to a spot
#    further down - see step 3).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Insert the new sentence about PlayerField / Goalkeeper right
#    after "...se desplazan. " and before the "E" that begins
#    "Existen varias instancias de Player".
# ------------------------------------------------------------------
$marker = "chequeo de límites del campo por donde se desplazan. E"
$full = $d.Content.Text
$markerIdx = $full.IndexOf($marker)
if ($markerIdx -lt 0) {
    throw "Could not locate insertion marker for the new sentence."
}
$insPoint = $markerIdx + ($marker.Length - 1)   # right before the trailing "E"

$ins = $d.Range($insPoint, $insPoint)
$ins.InsertBefore("De la clase ")

$full = $d.Content.Text
$p1 = $full.IndexOf("De la clase ") + "De la clase ".Length
$r = $d.Range($p1, $p1)
$r.InsertBefore("Player")
$d.Range($p1, $p1 + "Player".Length).Font.Italic = 1

$full = $d.Content.Text
$p2 = $full.IndexOf("Player heredan")
if ($p2 -lt 0) {
    $p2 = $p1 + "Player".Length
} else {
    $p2 = $p2 + "Player".Length
}
$r = $d.Range($p2, $p2)
$r.InsertBefore(" heredan ")

$full = $d.Content.Text
$p3 = $full.IndexOf(" heredan ", $p2) + " heredan ".Length
$r = $d.Range($p3, $p3)
$r.InsertBefore("PlayerField")
$d.Range($p3, $p3 + "PlayerField".Length).Font.Italic = 1

$full = $d.Content.Text
$p4 = $p3 + "PlayerField".Length
$r = $d.Range($p4, $p4)
$r.InsertBefore(" y ")

$full = $d.Content.Text
$p5 = $full.IndexOf(" y ", $p4) + " y ".Length
$r = $d.Range($p5, $p5)
$r.InsertBefore("Goalkeeper")
$d.Range($p5, $p5 + "Goalkeeper".Length).Font.Italic = 1

$full = $d.Content.Text
$p6 = $p5 + "Goalkeeper".Length
$r = $d.Range($p6, $p6)
$tailText = " que poseen diferentes límites de movimiento dentro de la cancha ya que el arquero solo puede moverse dentro del área grande"
$r.InsertBefore($tailText)

# ------------------------------------------------------------------
# 3) Drop a fresh "_GoBack" bookmark (zero-length) right after
#    "...área grande" - this is where the edit above left off.
# ------------------------------------------------------------------
$full = $d.Content.Text
$p7 = $p6 + $tailText.Length
$bmRange = $d.Range($p7, $p7)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4) Finish the sentence with ". " before the remaining "Existen..."
# ------------------------------------------------------------------
$r = $d.Range($p7, $p7)
$r.InsertBefore(". ")
